# Refactor synthetic "statut" array (3-color set): rename the "black"
# color entry (emoji + label) to a "blue" one.
#   ⬛ (black square)  -> 📘 (blue book)     statut_label "noir" -> "bleu"
#   🟥 (red square)    -> 📕 (red book)      statut_label "rouge" (unchanged)
#   🟩 (green square)  -> 📗 (green book)    statut_label "vert"  (unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "statut" (column A) currently holds the black-square emoji,
# and whose "statut_label" (column B) currently holds "noir".
$noirRows = @(2, 3, 6, 7, 9, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $noirRows) {
    $ws.Cells.Item($r, 1).Value = "📘"
    $ws.Cells.Item($r, 2).Value = "bleu"
}

# Rows whose "statut" (column A) holds the red-square emoji -> red book.
$rougeRows = @(4, 5)
foreach ($r in $rougeRows) {
    $ws.Cells.Item($r, 1).Value = "📕"
}

# Rows whose "statut" (column A) holds the green-square emoji -> green book.
$vertRows = @(8)
foreach ($r in $vertRows) {
    $ws.Cells.Item($r, 1).Value = "📗"
}
